# Update "想去人数" (want-to-go count) figures for newly scraped data.
# Workbook has 4 sheets: 展览, 演出, 本地生活, 全部类型 (本地生活 unchanged).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 1206
$ws.Range("F4").Value  = 1272
$ws.Range("F6").Value  = 178
$ws.Range("F8").Value  = 10
$ws.Range("F9").Value  = 339
$ws.Range("F11").Value = 1266
$ws.Range("F12").Value = 28978
$ws.Range("F13").Value = 3932
$ws.Range("F16").Value = 483
$ws.Range("F20").Value = 24
$ws.Range("F27").Value = 62
$ws.Range("F29").Value = 663
$ws.Range("F31").Value = 103
$ws.Range("F32").Value = 540
$ws.Range("F33").Value = 77
$ws.Range("F35").Value = 633
$ws.Range("F38").Value = 7

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value  = 876
$ws.Range("F17").Value = 46
$ws.Range("F23").Value = 4249

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F9").Value  = 876
$ws.Range("F10").Value = 1206
$ws.Range("F11").Value = 1272
$ws.Range("F12").Value = 178
$ws.Range("F14").Value = 10
$ws.Range("F15").Value = 339
$ws.Range("F17").Value = 1266
$ws.Range("F26").Value = 46
$ws.Range("F27").Value = 46
$ws.Range("F28").Value = 483
$ws.Range("F31").Value = 24
$ws.Range("F38").Value = 62
$ws.Range("F39").Value = 663
$ws.Range("F42").Value = 103
$ws.Range("F45").Value = 77
$ws.Range("F47").Value = 633
